# Update the p-values in the second table (Fig 8 p-values) of the document.
# Each change targets the Range of a specific table cell and uses
# Find/Replace scoped to that cell, so only the numeric text is replaced
# (avoiding the cell-end marker and any ambiguity with duplicate values
# elsewhere in the document).

$d = $word.ActiveDocument

$table = $d.Tables.Item(2)

function Replace-CellText($tbl, $row, $col, $oldText, $newText) {
    $cellRange = $tbl.Cell($row, $col).Range
    # Replace = 1 (wdReplaceOne) so the match is confined to this single
    # cell's Range instead of continuing into the rest of the document
    # (wdReplaceAll can keep matching past the end of the supplied range).
    $cellRange.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 0, $false, $newText, 1) | Out-Null
}

# Row 2 (Compilance (scores)), column 4 (socioeconomic): 0.40 -> 0.27
Replace-CellText $table 2 4 "0.40" "0.27"

# Row 3 (GDP 2016), column 4 (socioeconomic): 0.38 -> 0.87
Replace-CellText $table 3 4 "0.38" "0.87"

# Row 4 (OHI 2016), column 4 (socioeconomic): 0.99 -> 0.14
Replace-CellText $table 4 4 "0.99" "0.14"

# Row 5 (OHI economic 2016), column 4 (socioeconomic): 0.25 -> 0.49
Replace-CellText $table 5 4 "0.25" "0.49"

# Row 6 (Readiness), column 4 (socioeconomic): 0.36 -> 0.37
Replace-CellText $table 6 4 "0.36" "0.37"

# Row 7 (Technical Development), column 3 (institutional): 0.08 -> 0.06
Replace-CellText $table 7 3 "0.08" "0.06"

# Row 8 (Vulnerability), column 4 (socioeconomic): 0.23 -> 0.05
Replace-CellText $table 8 4 "0.23" "0.05"
